$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '65.585.71'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +2.84%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.214.40'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +2.67%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '600.99'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.50%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '152.80'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +4.83%  '

$ws.Range('E7').Value = '  -0.01%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.207.84'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +2.64%  '

$ws.Range('E9').Value = '  +3.14%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.168'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +3.88%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.19'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +8.59%  '

$ws.Range('E12').Value = '  +3.27%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000255'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +3.08%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '39.51'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +6.94%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.747.60'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +2.67%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.326.90'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +6.12%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '7.42'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +4.68%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '65.198.59'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +2.58%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '484.21'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +4.38%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '15.03'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +5.48%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.772'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +6.22%  '

$ws.Range('E23').Value = '  +6.46%  '

$ws.Range('E24').Value = '  +13.89%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '13.78'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +5.96%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '83.70'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.98%  '

$ws.Range('E27').Value = '  +10.72%  '

$ws.Range('E28').Value = '  +0.47%  '

$ws.Range('E29').Value = '  +4.55%  '

$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.55'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +7.94%  '

$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.29'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +3.80%  '

$ws.Range('E32').Value = '  +0.05%  '

$ws.Range('E33').Value = '  +10.57%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '28.85'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +7.30%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0₃0889'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.66%  '

$ws.Range('E36').Value = '  +5.06%  '

$ws.Range('B37').Value = 'dogwifhat'
$ws.Range('C37').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.52'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +3.64%  '

$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.41'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +4.53%  '

$ws.Range('B39').Value = 'Filecoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.37'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +6.09%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '485.47'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +10.04%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '52.34'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +4.06%  '

$ws.Range('E42').Value = '  +8.89%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.305'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +11.29%  '

$ws.Range('E44').Value = '  +4.24%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.952.58'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.49%  '

$ws.Range('E46').Value = '  +4.03%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '39.38'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +8.87%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '132.13'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +5.23%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.34'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +7.80%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '25.48'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +4.90%  '
